# This script applies targeted updates to column E ("Moyenne de l'etudiant")
# of the "BDA" worksheet, matching the corrective data fix described in the
# commit message: correctif probleme insertion matiere dans bd et mise a
# jour semestre etudiant lors de l'inscription.
#
# It rewrites the student average ("Moyenne") values for rows 4 through 63
# (skipping rows 35 and 39 which are unchanged) to their corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDA")

$updates = @{
    4 = 9;
    5 = 20;
    6 = 18;
    7 = 6;
    8 = 8;
    9 = 12;
    10 = 18;
    11 = 9;
    12 = 18;
    13 = 11;
    14 = 5;
    15 = 5;
    16 = 7;
    17 = 7;
    18 = 12;
    19 = 16;
    20 = 19;
    21 = 9;
    22 = 5;
    23 = 19;
    24 = 8;
    25 = 13;
    26 = 7;
    27 = 20;
    28 = 18;
    29 = 9;
    30 = 5;
    31 = 18;
    32 = 18;
    33 = 19;
    34 = 20;
    36 = 16;
    37 = 5;
    38 = 8;
    40 = 20;
    41 = 18;
    42 = 9;
    43 = 19;
    44 = 19;
    45 = 14;
    46 = 18;
    47 = 17;
    48 = 18;
    49 = 5;
    50 = 6;
    51 = 6;
    52 = 6;
    53 = 7;
    54 = 7;
    55 = 14;
    56 = 12;
    57 = 20;
    58 = 17;
    59 = 15;
    60 = 17;
    61 = 6;
    62 = 13;
    63 = 11;
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}
